$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 561
$ws.Range("A561").Value = 560
$ws.Range("B561").Value = 'Sunday, Jan 15'
$ws.Range("C561").Value = '4:10 PM'
$ws.Range("D561").Value = 'FR3365'
$ws.Range("E561").Value = 'Berlin'
$ws.Range("F561").Value = '(BER)'
$ws.Range("G561").Value = 'Ryanair '
$ws.Range("H561").Value = 'B738'
$ws.Range("I561").Value = '(9H-QES)'
$ws.Range("J561").Value = '4:08 PM'
$ws.Range("L561").Value = '0 hours, -2 minutes'
$ws.Range("K561").Borders.LineStyle = -4142
$ws.Range("M561").Borders.LineStyle = -4142

# Row 562
$ws.Range("A562").Value = 561
$ws.Range("B562").Value = 'Sunday, Jan 15'
$ws.Range("C562").Value = '4:30 PM'
$ws.Range("D562").Value = 'FR1643'
$ws.Range("E562").Value = 'Vienna'
$ws.Range("F562").Value = '(VIE)'
$ws.Range("G562").Value = 'Ryanair '
$ws.Range("H562").Value = 'B738'
$ws.Range("I562").Value = '(SP-RST)'
$ws.Range("J562").Value = '4:55 PM'
$ws.Range("L562").Value = '0 hours, 25 minutes'
$ws.Range("K562").Borders.LineStyle = -4142
$ws.Range("M562").Borders.LineStyle = -4142

# Row 563
$ws.Range("A563").Value = 562
$ws.Range("B563").Value = 'Sunday, Jan 15'
$ws.Range("C563").Value = '4:30 PM'
$ws.Range("D563").Value = 'FR6232'
$ws.Range("E563").Value = 'Palermo'
$ws.Range("F563").Value = '(PMO)'
$ws.Range("G563").Value = 'Buzz '
$ws.Range("H563").Value = 'B38M'
$ws.Range("I563").Value = '(SP-RZH)'
$ws.Range("J563").Value = '4:44 PM'
$ws.Range("L563").Value = '0 hours, 14 minutes'
$ws.Range("K563").Borders.LineStyle = -4142
$ws.Range("M563").Borders.LineStyle = -4142

# Row 564
$ws.Range("A564").Value = 563
$ws.Range("B564").Value = 'Sunday, Jan 15'
$ws.Range("C564").Value = '4:45 PM'
$ws.Range("D564").Value = 'LO3507'
$ws.Range("E564").Value = 'Olsztyn'
$ws.Range("F564").Value = '(SZY)'
$ws.Range("G564").Value = 'LOT '
$ws.Range("H564").Value = 'E75S'
$ws.Range("I564").Value = '(SP-LIA)'
$ws.Range("J564").Value = '4:32 PM'
$ws.Range("L564").Value = '0 hours, -13 minutes'
$ws.Range("K564").Borders.LineStyle = -4142
$ws.Range("M564").Borders.LineStyle = -4142

# Row 565
$ws.Range("A565").Value = 564
$ws.Range("B565").Value = 'Sunday, Jan 15'
$ws.Range("C565").Value = '4:45 PM'
$ws.Range("D565").Value = 'U22114'
$ws.Range("E565").Value = 'London'
$ws.Range("F565").Value = '(LTN)'
$ws.Range("G565").Value = 'easyJet '
$ws.Range("H565").Value = 'A320'
$ws.Range("I565").Value = '(G-EZWX)'
$ws.Range("J565").Value = '4:47 PM'
$ws.Range("L565").Value = '0 hours, 2 minutes'
$ws.Range("K565").Borders.LineStyle = -4142
$ws.Range("M565").Borders.LineStyle = -4142

# Row 566
$ws.Range("A566").Value = 565
$ws.Range("B566").Value = 'Sunday, Jan 15'
$ws.Range("C566").Value = '4:50 PM'
$ws.Range("D566").Value = 'KL1996'
$ws.Range("E566").Value = 'Amsterdam'
$ws.Range("F566").Value = '(AMS)'
$ws.Range("G566").Value = 'KLM '
$ws.Range("H566").Value = 'E190'
$ws.Range("I566").Value = '(PH-EXC)'
$ws.Range("J566").Value = '4:52 PM'
$ws.Range("L566").Value = '0 hours, 2 minutes'
$ws.Range("K566").Borders.LineStyle = -4142
$ws.Range("M566").Borders.LineStyle = -4142

# Row 567
$ws.Range("A567").Value = 566
$ws.Range("B567").Value = 'Sunday, Jan 15'
$ws.Range("C567").Value = '5:05 PM'
$ws.Range("D567").Value = 'FR1812'
$ws.Range("E567").Value = 'London'
$ws.Range("F567").Value = '(LTN)'
$ws.Range("G567").Value = 'Ryanair '
$ws.Range("H567").Value = 'B738'
$ws.Range("I567").Value = '(SP-RKB)'
$ws.Range("J567").Value = '5:17 PM'
$ws.Range("L567").Value = '0 hours, 12 minutes'
$ws.Range("K567").Borders.LineStyle = -4142
$ws.Range("M567").Borders.LineStyle = -4142

# Row 568
$ws.Range("A568").Value = 567
$ws.Range("B568").Value = 'Sunday, Jan 15'
$ws.Range("C568").Value = '5:10 PM'
$ws.Range("D568").Value = 'FR7115'
$ws.Range("E568").Value = 'Malta'
$ws.Range("F568").Value = '(MLA)'
$ws.Range("G568").Value = 'Ryanair '
$ws.Range("H568").Value = 'B38M'
$ws.Range("I568").Value = '(9H-VUU)'
$ws.Range("J568").Value = '5:57 PM'
$ws.Range("L568").Value = '0 hours, 47 minutes'
$ws.Range("K568").Borders.LineStyle = -4142
$ws.Range("M568").Borders.LineStyle = -4142

# Row 569
$ws.Range("A569").Value = 568
$ws.Range("B569").Value = 'Sunday, Jan 15'
$ws.Range("C569").Value = '5:25 PM'
$ws.Range("D569").Value = 'W65017'
$ws.Range("E569").Value = 'Birmingham'
$ws.Range("F569").Value = '(BHX)'
$ws.Range("G569").Value = 'Wizz Air '
$ws.Range("H569").Value = 'A21N'
$ws.Range("I569").Value = '(HA-LVG)'
$ws.Range("J569").Value = '5:53 PM'
$ws.Range("L569").Value = '0 hours, 28 minutes'
$ws.Range("K569").Borders.LineStyle = -4142
$ws.Range("M569").Borders.LineStyle = -4142

# Row 570
$ws.Range("A570").Value = 569
$ws.Range("B570").Value = 'Sunday, Jan 15'
$ws.Range("C570").Value = '5:30 PM'
$ws.Range("D570").Value = 'FR7660'
$ws.Range("E570").Value = 'Newcastle'
$ws.Range("F570").Value = '(NCL)'
$ws.Range("G570").Value = 'Buzz '
$ws.Range("H570").Value = 'B38M'
$ws.Range("I570").Value = '(SP-RZB)'
$ws.Range("J570").Value = '5:35 PM'
$ws.Range("L570").Value = '0 hours, 5 minutes'
$ws.Range("K570").Borders.LineStyle = -4142
$ws.Range("M570").Borders.LineStyle = -4142

# Row 571
$ws.Range("A571").Value = 570
$ws.Range("B571").Value = 'Sunday, Jan 15'
$ws.Range("C571").Value = '5:45 PM'
$ws.Range("D571").Value = 'FR3054'
$ws.Range("E571").Value = 'Barcelona'
$ws.Range("F571").Value = '(BCN)'
$ws.Range("G571").Value = 'Ryanair '
$ws.Range("H571").Value = 'B738'
$ws.Range("I571").Value = '(EI-DPG)'
$ws.Range("J571").Value = '5:46 PM'
$ws.Range("L571").Value = '0 hours, 1 minutes'
$ws.Range("K571").Borders.LineStyle = -4142
$ws.Range("M571").Borders.LineStyle = -4142

# Row 572
$ws.Range("A572").Value = 571
$ws.Range("B572").Value = 'Sunday, Jan 15'
$ws.Range("C572").Value = '5:45 PM'
$ws.Range("D572").Value = 'W65003'
$ws.Range("E572").Value = 'London'
$ws.Range("F572").Value = '(LTN)'
$ws.Range("G572").Value = 'Wizz Air '
$ws.Range("H572").Value = 'A21N'
$ws.Range("I572").Value = '(HA-LZI)'
$ws.Range("J572").Value = '5:51 PM'
$ws.Range("L572").Value = '0 hours, 6 minutes'
$ws.Range("K572").Borders.LineStyle = -4142
$ws.Range("M572").Borders.LineStyle = -4142
